$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.137.05"
$ws.Range("E2").Value = "  -4.44%  "
$ws.Range("D3").Value = "1.833.56"
$ws.Range("E3").Value = "  -3.01%  "
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").Value = "'329.39"
$ws.Range("E5").Value = "  -2.67%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "'0.4654"
$ws.Range("E7").Value = "  -1.97%  "
$ws.Range("D8").Value = "'0.3875"
$ws.Range("E8").Value = "  -2.99%  "
$ws.Range("D9").Value = "'46.23"
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("D10").Value = "'0.07897"
$ws.Range("E10").Value = "  -1.50%  "
$ws.Range("D11").Value = "'0.9624"
$ws.Range("E11").Value = "  -2.68%  "
$ws.Range("E12").Value = "  -4.70%  "
$ws.Range("D13").Value = "1.848.57"
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("D14").Value = "'5.673"
$ws.Range("E14").Value = "  -4.33%  "
$ws.Range("D15").Value = "'6.911"
$ws.Range("E15").Value = "  -2.43%  "
$ws.Range("D16").Value = "'0.06863"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").Value = "'86.76"
$ws.Range("E18").Value = "  -2.51%  "
$ws.Range("D19").Value = "'0.000009990"
$ws.Range("E19").Value = "  -1.89%  "
$ws.Range("E20").Value = "  -3.52%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "28.151.77"
$ws.Range("E22").Value = "  -4.45%  "
$ws.Range("D23").Value = "'5.337"
$ws.Range("E23").Value = "  -2.92%  "
$ws.Range("D24").Value = "'11.02"
$ws.Range("E24").Value = "  -5.26%  "
$ws.Range("D25").Value = "'2.095"
$ws.Range("E25").Value = "  -2.47%  "
$ws.Range("D26").Value = "2.045.31"
$ws.Range("E26").Value = "  -2.94%  "
$ws.Range("D27").Value = "'152.77"
$ws.Range("E27").Value = "  -2.93%  "
$ws.Range("D28").Value = "'19.27"
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D29").Value = "'5.772"
$ws.Range("E29").Value = "  -10.97%  "
$ws.Range("D30").Value = "'1.978"
$ws.Range("E30").Value = "  -3.26%  "
$ws.Range("D31").Value = "'117.36"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.09285"
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.9373"
$ws.Range("E33").Value = "  -5.89%  "
$ws.Range("D34").Value = "'5.309"
$ws.Range("E34").Value = "  -2.87%  "
$ws.Range("E35").Value = "  -4.32%  "
$ws.Range("D36").Value = "'3.358"
$ws.Range("E36").Value = "  -4.85%  "
$ws.Range("D37").Value = "'0.05944"
$ws.Range("E37").Value = "  -7.32%  "
$ws.Range("D38").Value = "'0.02152"
$ws.Range("E38").Value = "  -3.97%  "
$ws.Range("D39").Value = "'1.152"
$ws.Range("E39").Value = "  -3.89%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'7.685"
$ws.Range("E40").Value = "  -0.42%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5602"
$ws.Range("E41").Value = "  -3.75%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'9.940"
$ws.Range("E42").Value = "  -5.37%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1774"
$ws.Range("E43").Value = "  -2.43%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'1.236"
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'2.216"
$ws.Range("E45").Value = "  -8.10%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'11.62"
$ws.Range("E46").Value = "  -4.24%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5281"
$ws.Range("E47").Value = "  -3.79%  "
$ws.Range("D48").Value = "'0.07055"
$ws.Range("E48").Value = "  -3.79%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.839"
$ws.Range("E49").Value = "  -5.67%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'111.99"
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Value = "'1.001"
$ws.Range("E51").Value = "  -0.38%  "
